$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "29.196.67"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "1.848.04"
$ws.Range("E3").Value = "  -0.36%  "
Set-TextValue "D5" "245.63"
$ws.Range("E5").Value = "  +1.78%  "
Set-TextValue "D6" "0.7017"
$ws.Range("E6").Value = "  +0.61%  "
Set-TextValue "D7" "1.001"
$ws.Range("E7").Value = "  +0.06%  "
Set-TextValue "D8" "0.07718"
$ws.Range("E8").Value = "  +0.02%  "
Set-TextValue "D9" "0.3069"
$ws.Range("E9").Value = "  -0.08%  "
Set-TextValue "D10" "23.64"
$ws.Range("E10").Value = "  -0.45%  "
Set-TextValue "D11" "0.07815"
$ws.Range("E11").Value = "  +0.07%  "
Set-TextValue "D12" "92.95"
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.849.93"
$ws.Range("E13").Value = "  -0.34%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D14" "5.140"
$ws.Range("E14").Value = "  +0.86%  "
Set-TextValue "D15" "0.6862"
$ws.Range("E15").Value = "  +0.03%  "
Set-TextValue "D16" "6.601"
$ws.Range("E16").Value = "  +1.12%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "29.200.57"
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D18" "0.000008307"
$ws.Range("E18").Value = "  -1.15%  "
Set-TextValue "D19" "241.89"
$ws.Range("E19").Value = "  -2.91%  "
$ws.Range("D20").Value = "2.091.86"
$ws.Range("E20").Value = "  -1.12%  "
Set-TextValue "D21" "12.73"
$ws.Range("E21").Value = "  -0.49%  "
Set-TextValue "D22" "1.000"
Set-TextValue "D23" "7.523"
$ws.Range("E23").Value = "  +0.18%  "
Set-TextValue "D24" "1.000"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("E25").Value = "  +0.73%  "
Set-TextValue "D26" "159.29"
$ws.Range("E26").Value = "  -0.87%  "
Set-TextValue "D27" "8.825"
$ws.Range("E27").Value = "  -0.22%  "
Set-TextValue "D28" "18.30"
$ws.Range("E28").Value = "  -0.83%  "
Set-TextValue "D29" "1.537"
$ws.Range("E29").Value = "  -1.42%  "
Set-TextValue "D30" "4.226"
$ws.Range("E30").Value = "  -0.31%  "
Set-TextValue "D31" "4.183"
$ws.Range("E32").Value = "  +0.99%  "
Set-TextValue "D33" "0.05123"
$ws.Range("E33").Value = "  -1.52%  "
Set-TextValue "D34" "0.7922"
$ws.Range("E34").Value = "  +4.28%  "
Set-TextValue "D35" "1.910"
$ws.Range("E35").Value = "  +3.84%  "
$ws.Range("E36").Value = "  -1.58%  "
Set-TextValue "D37" "2.698"
$ws.Range("E37").Value = "  -0.52%  "
$ws.Range("D38").Value = "1.322.89"
$ws.Range("E38").Value = "  +8.97%  "
Set-TextValue "D39" "0.01871"
$ws.Range("E39").Value = "  +0.58%  "
Set-TextValue "D40" "2.714"
$ws.Range("E40").Value = "  -0.25%  "
Set-TextValue "D41" "0.9582"
$ws.Range("E41").Value = "  +6.88%  "
Set-TextValue "D42" "6.070"
$ws.Range("E42").Value = "  +10.22%  "
Set-TextValue "D43" "107.34"
$ws.Range("E43").Value = "  -2.33%  "
$ws.Range("E44").Value = "  +0.11%  "
Set-TextValue "D45" "9.708"
$ws.Range("E45").Value = "  +1.94%  "
$ws.Range("E46").Value = "  -1.00%  "
$ws.Range("D47").Value = "1.992.71"
$ws.Range("E47").Value = "  -1.00%  "
Set-TextValue "D48" "0.5182"
$ws.Range("E48").Value = "  +0.04%  "
Set-TextValue "D49" "64.35"
$ws.Range("E49").Value = "  -1.13%  "
Set-TextValue "D50" "1.762"
$ws.Range("E50").Value = "  +0.75%  "
Set-TextValue "D51" "6.990"
$ws.Range("E51").Value = "  -0.16%  "
